# Apply the product-backlog edits: update story names/descriptions, re-scale
# estimation values to T-shirt sizes, add start/finish dates, adjust wrap/row/
# column formatting and selection to match the authored workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Story Name (row 3 renamed) ---
$ws.Range("B3").Value = 'Query Time-specific Expenses via Natural Language'

# --- Descriptions (now multi-line; turn on wrap so the row can grow) ---
$ws.Range("C2").Value = 'As a user, 
I want to query my current account balance using natural language , 
so that I can quickly understand my financial situation.'
$ws.Range("C3").Value = 'As a user,
I would like to use natural language to look up my total expenses and specific expense categories for a particular period of time, 
so that I can manage my finances better.'
$ws.Range("C4").Value = 'As a user, 
I want the AI assistant to analyze my spending patterns and provide saving recommendations via natural language,
 so that I can make smarter financial decisions.'
$ws.Range("C5").Value = 'As a user, 
I want to record new transactions quickly via speech or text input, 
so that I don''t need to enter transactions manually.'
$ws.Range("C2:C5").WrapText = $true

# --- Estimation: numeric story points -> t-shirt sizes ---
$ws.Range("G2").Value = "small"
$ws.Range("G3").Value = "medium"
$ws.Range("G4").Value = "extra large"
$ws.Range("G5").Value = "medium"

# --- Date started (actual date): write as literal text, not an Excel date ---
$ws.Range("I2").Formula = "=""2025.3.13"""
$ws.Range("I3").Formula = "=""2025.3.13"""
$ws.Range("I4").Formula = "=""2025.3.14"""
$ws.Range("I5").Formula = "=""2025.3.15"""
$ws.Range("I2:I5").Copy()
$ws.Range("I2:I5").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Row height for the now-taller Story 2 description, column E width, selection ---
$ws.Rows("3:3").RowHeight = 100.8
$ws.Columns("E:E").ColumnWidth = 27.727120535714285
$ws.Range("E4").Select()

